$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) currently reads, through column M:
#   ... J1=Nombre Iz  K1=Puesto Iz  L1=Nombre D  M1=Puesto D
# New "Abajo Iz" / "Abajo D" labels are being added to the right-hand
# signature block, giving:
#   ... J1=Nombre Iz K1=Puesto Iz L1=Abajo Iz M1=Nombre D N1=Puesto D O1=Abajo D

# Drop the two existing trailing header cells (L1:M1) so they can be
# rebuilt in the right order, then insert 4 fresh columns in their place
# (L, M, N, O) for the new 4-column header tail.
$ws.Range("L1:M1").EntireColumn.Delete()
$ws.Range("L1:M1").EntireColumn.Insert()

$ws.Range("L1").Value = "Abajo Iz"
$ws.Range("M1").Value = "Nombre D"

# N1 / O1 are brand new cells beyond the old A1:M2 dimension, so they pick
# up the sheet's default (unstyled) formatting rather than the bold header
# style used by A1:M1.
$ws.Range("N1").Value = "Puesto D"
$ws.Range("O1").Value = "Abajo D"

# Restore the selection marker to match the post-edit workbook (it was
# parked at M3, now at M15).
$ws.Range("M15").Select()
